$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = -1
$ws.Range("G17").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("G24").Value = 0

$ws.Range("E17").Select()
